# Updated cryptos list on Mon Sep 23 23:44:24 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.366.92"
$ws.Range("E2").Value = "  -0.34%  "

$ws.Range("D3").Value = "2.648.50"
$ws.Range("E3").Value = "  +2.55%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'604.31"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.92%  "

$ws.Range("D6").Value = "'144.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.02%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'0.586"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.37%  "

$ws.Range("D9").Value = "2.647.87"
$ws.Range("E9").Value = "  +2.55%  "

$ws.Range("E10").Value = "  +1.78%  "

$ws.Range("E11").Value = "  +0.45%  "

$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("E13").Value = "  +3.41%  "

$ws.Range("D14").Value = "'27.40"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.29%  "

$ws.Range("D15").Value = "3.123.79"
$ws.Range("E15").Value = "  +2.58%  "

$ws.Range("D16").Value = "63.198.46"
$ws.Range("E16").Value = "  -0.44%  "

$ws.Range("E17").Value = "  -0.57%  "

$ws.Range("D18").Value = "2.650.91"
$ws.Range("E18").Value = "  +2.92%  "

$ws.Range("D19").Value = "'11.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.55%  "

$ws.Range("D20").Value = "'4.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.20%  "

$ws.Range("D21").Value = "'342.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("D22").Value = "'6.84"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.00%  "

$ws.Range("D23").Value = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.20%  "

$ws.Range("D24").Value = "'67.16"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.02%  "

$ws.Range("D25").Value = "'1.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.82%  "

$ws.Range("E26").Value = "  -2.64%  "

$ws.Range("D27").Value = "'8.68"
$ws.Range("D27").Style = "Normal"

$ws.Range("E28").Value = "  -0.05%  "

$ws.Range("D29").Value = "'543.48"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +15.23%  "

$ws.Range("E30").Value = "  +0.10%  "

$ws.Range("D31").Value = "'7.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.38%  "

$ws.Range("E32").Value = "  +4.17%  "

$ws.Range("E33").Value = "  +7.64%  "

$ws.Range("D34").Value = "0.0₃0811"
$ws.Range("E34").Value = "  +1.13%  "

$ws.Range("D35").Value = "'171.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.69%  "

$ws.Range("D36").Value = "'5.12"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +13.01%  "

$ws.Range("E37").Value = "  +1.01%  "

$ws.Range("D38").Value = "'1.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "'19.13"
$ws.Range("D39").Style = "Normal"

$ws.Range("E40").Value = "  +6.40%  "

$ws.Range("D41").Value = "'172.40"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.79%  "

$ws.Range("D43").Value = "'3.77"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.39%  "

$ws.Range("D44").Value = "'22.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.58%  "

$ws.Range("D45").Value = "'0.0579"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.71%  "

$ws.Range("D46").Value = "'0.631"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.12%  "

$ws.Range("D47").Value = "'0.0962"
$ws.Range("D47").Style = "Normal"

$ws.Range("D48").Value = "'0.0240"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.82%  "

$ws.Range("D49").Value = "'18.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.96%  "

$ws.Range("E50").Value = "  +2.67%  "

$ws.Range("E51").Value = "  -0.93%  "
